$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 31
$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "RJ TV 2"
$ws.Cells.Item($row, 3).Value = "Trânsito"
$ws.Cells.Item($row, 4).Value = "2025-04-01T18:12"
$ws.Cells.Item($row, 5).Value = "Neutro"
$ws.Cells.Item($row, 6).Value = "Caminhão invade calçadão de Campos. Motorista teria errado o caminho. Comerciantes ficaram sem energia elétrica. Repórter *ao vivo*. Carreta foi retirada e energia foi restabelecida. Um poste caiu. Muitos comerciantes foram prejudicados hoje. Motorista é do Rio Grande do Sul e GPS indicou caminho errado. Por não poder pegar a Av. Arthur Bernardes, GPS indicou outro caminho. Entrevista com pessoas que passavam pelo local e com presidente da Acic, Maurício Cabral. Motorista foi levado à delegacia. Guarda Municipal informou que motorista foi notificado por transitar no passeio em local não permitido.  "
